# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.906.09"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.787.41"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D5").Value = "'357.81"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'109.51"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("D7").Value = "'0.562"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "'0.0846"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "'19.49"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "'7.59"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").Value = "3.229.08"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "2.791.56"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "'0.951"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "51.860.16"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").Value = "'13.15"
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'270.74"
$ws.Range("D24").Value = "'70.17"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'2.76"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").Value = "'26.44"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'0.166"
$ws.Range("E28").Value = "  +18.53%  "
$ws.Range("D29").Value = "'10.30"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'52.03"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").Value = "'34.75"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("E36").Value = "  -5.44%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'18.75"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E40").Value = "  -4.51%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "'119.45"
$ws.Range("E44").Value = "  -4.62%  "
$ws.Range("E45").Value = "  -6.51%  "
$ws.Range("D46").Value = "2.078.20"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'3.27"
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("D48").Value = "'2.23"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("D50").Value = "'0.949"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("D51").Value = "'8.68"
$ws.Range("E51").Value = "  -4.67%  "
